$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append a new row of mail-log data ---
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A15").Value = "Algemene vraag"
$logs.Range("B15").Value = "klachten@testbedrijf123.nl"
$logs.Range("D15").Value = "Intern verzoek / Actie voor medewerker"
$logs.Range("F15").Value = "2025-08-19 19:58:52"
$logs.Range("G15").Value = "Nee"
$logs.Range("H15").Value = "Ja"
$logs.Range("I15").Value = "Nee"
$logs.Range("J15").Value = "Nee"

# --- Extend the existing conditional-formatting ranges to cover the new row ---
$columns = "D", "G", "H", "I", "J"
foreach ($col in $columns) {
    $oldRange = $col + "2:" + $col + "14"
    $newRange = $col + "2:" + $col + "15"
    $fcs = $logs.Range($oldRange).FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($logs.Range($newRange))
    }
}

# --- Sheet "Dashboard": bump the aggregate count for this category ---
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 14
